$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two stray C2/C3 values (naive component forecaster bug fix)
$ws.Range("C2").ClearContents()
$ws.Range("C3").ClearContents()

# Tiny floating-point corrections from the re-run of the forecaster
$ws.Range("C4").Value = -7.266312015249799
$ws.Range("E5").Value = 12.2165830639507
$ws.Range("C6").Value = 9.469137444079955
$ws.Range("E6").Value = 8.07926457985193
$ws.Range("C7").Value = 3.358206407534969
$ws.Range("C9").Value = 3.901355411819685
$ws.Range("E10").Value = 4.089819750351809
$ws.Range("E11").Value = 2.313009565865709
$ws.Range("C12").Value = 5.246209615995689
$ws.Range("C13").Value = 4.862559663742938
$ws.Range("E13").Value = 4.112897401876769
$ws.Range("C14").Value = 2.76474001115945
$ws.Range("E14").Value = 1.643374185611379
$ws.Range("C15").Value = -7.260793671746447
$ws.Range("E15").Value = -5.080986607234461
$ws.Range("C16").Value = 4.097586525396246
$ws.Range("C17").Value = 7.824284864703768
$ws.Range("C18").Value = -1.245022353133318
$ws.Range("E19").Value = 2.806286889124987
